$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.015390396118164
$ws.Range("B1").Value = 1.368349313735962
$ws.Range("C1").Value = 2.273638963699341
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.880639433860779
